# FOR V22 WARDEN LIGHT U10! Credits Menu, Changelog & Support!
#
# 1) The stray "_GoBack" bookmark that Word drops at the last edit position
#    (currently sitting between " Voxel Game " and "Light Version" in the
#    title) is removed from its old spot...
# 2) ...two new "No Spacing" paragraphs are appended after the Testers
#    credit line: a freesound.org / Erokia audio-credit line, and a blank
#    paragraph that becomes the new home of the "_GoBack" bookmark (i.e.
#    the bookmark just follows the last edit, as Word does automatically).

$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark from the title paragraph -------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append the audio-credit paragraph + the bookmark paragraph -------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParasXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
    "<w:r><w:t xml:space=`"preserve`">Audio material was provided from freesound.org by </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Erokia</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
    "<w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"

$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.InsertXML($newParasXml)
